$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")

# ---- Row 1 (message key names) ----
# A1, B1 already correct ("msg_pass_cart_badge", "msg_pass_cart_badge2" -> needs update)
$ws.Range("B1").Value = "msg_not_pass_cart_badge2"
$ws.Range("C1").Value = "msg_not_pass_display_sizeProduct"
$ws.Range("D1").Value = "msg_pass_display_sizeProduct"
$ws.Range("E1").Value = "msg_not_pass_display_nameProduct"
$ws.Range("F1").Value = "msg_pass_display_nameProduct"
$ws.Range("G1").Value = "msg_not_pass_display_priceProduct"
$ws.Range("H1").Value = "msg_pass_display_priceProduct"
$ws.Range("I1").Value = "msg_not_pass_display_quantityProduct"
$ws.Range("J1").Value = "msg_pass_display_quantityProduct"
$ws.Range("K1").Value = "msg_not_pass_display_btn"
$ws.Range("L1").Value = "msg_pass_display_btn"
$ws.Range("M1").Value = "msg_pass_click_btn"

# ---- Row 2 (message text values) ----
$ws.Range("B2").Value = "❌ The total number off products in the cart is incorrect!"
$ws.Range("C2").Value = "`tThe number of products displayed is incorrect; `n`t{0} products were added, but {1} are shown."
$ws.Range("D2").Value = "✅ Number of products on display: {0}"
$ws.Range("E2").Value = '\tProduct name number {0} is displayed incorrectly.\n\tAdded name: {1}.\n\tDisplayed name: {2}.'
$ws.Range("F2").Value = "Name Product {0}: {1}"
$ws.Range("G2").Value = "`tProduct price number {0} is displayed incorrectly.`n`tAdded price: {1}.`n`tDisplayed price: {2}."
$ws.Range("H2").Value = "Price Product {0}: {1}"
$ws.Range("I2").Value = "`tProduct quantity number {0} is displayed incorrectly.`n`tAdded quantity: {1}.`n`tDisplayed quantity: {2}."
$ws.Range("J2").Value = "Quantity Product: {1}"
$ws.Range("K2").Value = "❌ Button {0} is not displayed"
$ws.Range("L2").Value = "✅ Button {0} is displayed"
$ws.Range("M2").Value = "✅ Button {0} is clicked"

# ---- Formatting: new font (black Aptos Narrow) + wrap text on header cells F1:M1 ----
$ws.Range("A1:M1").WrapText = $true
$ws.Range("A2:M2").WrapText = $true
$ws.Range("F1:M1").Font.Color = 0

# ---- Row heights ----
$ws.Rows.Item(1).RowHeight = 68
$ws.Rows.Item(2).RowHeight = 187

# ---- Selection ----
$ws.Range("F2").Select() | Out-Null
